$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.318.05"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('D3').Value = "'3.896.21"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'531.18"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +9.38%  '
$ws.Range('D6').Value = "'144.57"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('D7').Value = "'0.612"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.57%  '
$ws.Range('E9').Value = '  -2.91%  '
$ws.Range('E10').Value = '  -2.33%  '
$ws.Range('D11').Value = "'0.0000333"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -5.98%  '
$ws.Range('D12').Value = "'42.09"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').Value = "'4.516.78"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D15').Value = "'3.898.85"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = "'13.97"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.87%  '
$ws.Range('E17').Value = '  +6.70%  '
$ws.Range('E18').Value = '  -1.44%  '
$ws.Range('D19').Value = "'19.76"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').Value = "'69.275.08"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('D21').Value = "'424.93"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('E22').Value = '  -5.40%  '
$ws.Range('E23').Value = '  -4.14%  '
$ws.Range('D24').Value = "'87.85"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.23%  '
$ws.Range('E25').Value = '  +8.85%  '
$ws.Range('D26').Value = "'11.48"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -7.67%  '
$ws.Range('D27').Value = "'10.54"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.15%  '
$ws.Range('E28').Value = '  -2.31%  '
$ws.Range('D29').Value = "'687.94"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.61%  '
$ws.Range('D30').Value = "'13.18"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('E32').Value = '  -2.90%  '
$ws.Range('D33').Value = "'68.55"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +10.94%  '
$ws.Range('D34').Value = "'0.0₃0863"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.93%  '
$ws.Range('E35').Value = '  +7.70%  '
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('E38').Value = '  +2.51%  '
$ws.Range('D39').Value = "'0.998"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').Value = "'3.29"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +7.55%  '
$ws.Range('D42').Value = "'0.0484"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.40%  '
$ws.Range('D43').Value = "'3.22"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +7.69%  '
$ws.Range('D44').Value = "'2.78"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -7.44%  '
$ws.Range('D45').Value = "'3.41"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.87%  '
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('D47').Value = "'0.000278"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +12.58%  '
$ws.Range('D48').Value = "'2.99"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +6.83%  '
$ws.Range('D49').Value = "'2.756.89"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +15.03%  '
$ws.Range('E50').Value = '  -5.74%  '
$ws.Range('D51').Value = "'144.80"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.40%  '
